# adj_pronoun sheet: re-point the dep_rel AutoFilter from "MOS" to "POSDEP"
# (correcting the role assigned to adj-concatenated pronouns whose dep role
# is NPOSTMOD / POSDEP), which also flips which rows the filter hides, and
# move the active selection to reflect the new first visible cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("dep_rel") is the 5th column of the filtered range A1:F413.
# Passing Criteria1 as an array with Operator = xlFilterValues (7) produces
# the discrete-values <filters><filter val="..."/></filters> form (as
# opposed to a <customFilters> comparison), matching a normal "check just
# this value" filter selection in the Excel UI.
$ws.Range("A1:F413").AutoFilter(5, @("POSDEP"), 7)

# Update the active cell/selection to F103.
$ws.Range("F103").Select()
